# This script applies a set of precise edits to the document, mirroring
# the target unified diff:
#  0. Register a new numbering definition (abstractNumId 1 / numId 2),
#     re-submitting the (unmodified) document body alongside it, since
#     the engine only honours an updated numbering.xml "side part" when
#     the Range passed to InsertXML is the full document Content range.
#  1. Paragraph 1: split " Définition CPU" run into 3 runs with spell-check
#     proofErr markers around "Définition".
#  2. Paragraph 3: split "Description protocole SPI" run into 3 runs with
#     spell-check proofErr markers around "protocole".
#  3. Paragraph 4: merge "GP" + bookmark + "IO" into a single run "GPIO"
#     (dropping the now-unused bookmark from this paragraph).
#  4. After paragraph 4, insert two empty paragraphs, an "Attente :"
#     paragraph, and a new bulleted list paragraph ("Information  Capteurs
#     pressions") that reuses the _GoBack bookmark.
#
# NOTE: steps 1-4 are deliberately applied as narrow-range InsertXML calls
# (covering only the run content of each paragraph, not the paragraph mark)
# so that paragraph-level attributes (w:rsidR, w:rsidP, w:pPr, ...) that the
# reference diff leaves untouched are preserved exactly. Reading back
# Range.WordOpenXML after edits normalizes/merges runs and drops proofErr
# markers, so it must not be used as an intermediate snapshot; that is why
# the numbering registration (step 0) happens first, before any of the
# surgical run-level edits.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 0: register abstractNumId=1 / numId=2 in numbering.xml
# ---------------------------------------------------------------------
$origBody = '<w:p w:rsidR="003F7452" w:rsidRDefault="00A62882" w:rsidP="00A62882"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> Définition CPU</w:t></w:r></w:p><w:p w:rsidR="005B45B8" w:rsidRDefault="005B45B8" w:rsidP="00A62882"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Description USB2.0</w:t></w:r></w:p><w:p w:rsidR="00B66A83" w:rsidRDefault="00B66A83" w:rsidP="00A62882"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Description protocole SPI</w:t></w:r></w:p><w:p w:rsidR="0074315D" w:rsidRDefault="0074315D" w:rsidP="00A62882"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Description entrées/sorties GP</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>IO</w:t></w:r></w:p><w:sectPr w:rsidR="0074315D"><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1417" w:right="1417" w:bottom="1417" w:left="1417" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr>'

$numAbs0 = '<w:abstractNum w:abstractNumId="0"><w:nsid w:val="14BF28FE"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="9F7AB046"/><w:lvl w:ilvl="0" w:tplc="44BC46D4"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1-"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="04090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%2."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="2" w:tplc="0409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%3."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="2160" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="3" w:tplc="0409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%4."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="4" w:tplc="04090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%5."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="5" w:tplc="0409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%6."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="4320" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="6" w:tplc="0409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%7."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="7" w:tplc="04090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%8."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="8" w:tplc="0409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%9."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="6480" w:hanging="180"/></w:pPr></w:lvl></w:abstractNum>'

$numAbs1 = '<w:abstractNum w:abstractNumId="1"><w:nsid w:val="629A05BD"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="276840F6"/><w:lvl w:ilvl="0" w:tplc="CDDCF016"><w:numFmt w:val="bullet"/><w:lvlText w:val="-"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Calibri" w:cs="Calibri" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="04090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="04090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0A7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="04090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0B7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="04090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="04090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0A7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="04090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0B7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="04090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="04090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0A7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum>'

$numNums = '<w:num w:numId="1"><w:abstractNumId w:val="0"/></w:num><w:num w:numId="2"><w:abstractNumId w:val="1"/></w:num>'

$numberingXml = '<w:numbering xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $numAbs0 + $numAbs1 + $numNums + '</w:numbering>'

$step0Pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $origBody + '</w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/numbering.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.numbering+xml"><pkg:xmlData>' + $numberingXml + '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($step0Pkg)

# ---------------------------------------------------------------------
# Common XML envelope used for the remaining surgical, narrow-range edits
# ---------------------------------------------------------------------
$pkgHead = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Step 1: Paragraph 1 -> " " / "Définition" (spellStart/spellEnd) / " CPU"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1Range = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$inner1 = '<w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Définition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> CPU</w:t></w:r></w:p>'
$xml1 = $pkgHead + $inner1 + $pkgTail
$p1Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# Step 2: Paragraph 3 -> "Description " / "protocole" (spellStart/spellEnd) / " SPI"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3Range = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$inner2 = '<w:p><w:r><w:t xml:space="preserve">Description </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>protocole</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> SPI</w:t></w:r></w:p>'
$xml2 = $pkgHead + $inner2 + $pkgTail
$p3Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# Step 3: Paragraph 4 -> merge "GP" + bookmark + "IO" into "GPIO"
#         (drop the now-unused bookmark from this paragraph)
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4Range = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$inner3 = '<w:p><w:r><w:t>Description entrées/sorties GPIO</w:t></w:r></w:p>'
$xml3 = $pkgHead + $inner3 + $pkgTail
$p4Range.InsertXML($xml3)

# ---------------------------------------------------------------------
# Step 4: Insert new paragraphs after paragraph 4 (before the section break)
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$insertPoint = $d.Range($p4.Range.End, $p4.Range.End)
$inner4 = '<w:p/><w:p/><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Attente</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Information  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Capteurs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>press</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ions</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$xml4 = $pkgHead + $inner4 + $pkgTail
$insertPoint.InsertXML($xml4)
